$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Delete the header row (row 1): "Departing From | Departing Month | Departing day | Returning from | Returning month | Returning Day"
# This shifts the two data rows up, so former row2 becomes row1 and former row3 becomes row2.
$ws.Rows.Item(1).Delete()

# Apply a Text number format to columns C and F (used for day-of-month values)
$ws.Columns.Item(3).NumberFormat = "@"
$ws.Columns.Item(6).NumberFormat = "@"

# Set page orientation to portrait
$ws.PageSetup.Orientation = 1

# Update selection to match the resulting view
$ws.Range("I25").Select()
